$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only codigo/produto changed
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "02285"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "INTENSE BAS LIQ MATE CAMUFL POP 320 20ml"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "5"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "3"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "H"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "02286"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "INTENSE BAS LIQ MATE CAMUFL POP 330 20ml"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "65"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "asa"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "6"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "A"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "02289"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "COMB INTENSE GLITTER CARNAVAL"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "85"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "psa"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "8"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "15"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "5"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "N"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "02483"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "MAKE"

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "16"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "6"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "N"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "02484"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "6"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "MAKE"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "10"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "17"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "7"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "f"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "02485"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "7"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "MAKE"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "8"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "18"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "8"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "G"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "02485"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "8"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "MAKE"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "19"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "9"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "h"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "02485"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "9"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "MAKE"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "8"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "20"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "10"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "l"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "02485"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "10"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "MAKE"

# Row 11
$ws.Range("A11").Value = 3
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "21"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "11"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "q"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "02485"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "INTENSE LAP P/OLHO POP PRETO 1,1g"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "11"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "MAKE"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "6"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "8"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "5"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "B"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "55559"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "HER CODE"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "12"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "PERFUME"
